$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 45-46: coin order swap (Decentraland <-> WEMIXTOKEN) + value updates
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"

# Column D (Price) updates - forced as Text to preserve exact formatting
# (trailing zeros, thousand-dot notation, precision) like the source data feed
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "27.822.26"
$cell.ClearFormats()
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "1.878.88"
$cell.ClearFormats()
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "1.009"
$cell.ClearFormats()
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "335.69"
$cell.ClearFormats()
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "1.008"
$cell.ClearFormats()
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4698"
$cell.ClearFormats()
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.3947"
$cell.ClearFormats()
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "45.65"
$cell.ClearFormats()
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.08023"
$cell.ClearFormats()
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "1.012"
$cell.ClearFormats()
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "22.05"
$cell.ClearFormats()
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "1.875.78"
$cell.ClearFormats()
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "6.023"
$cell.ClearFormats()
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "7.301"
$cell.ClearFormats()
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = "1.010"
$cell.ClearFormats()
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "89.08"
$cell.ClearFormats()
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "0.06741"
$cell.ClearFormats()
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "0.00001048"
$cell.ClearFormats()
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "17.39"
$cell.ClearFormats()
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "27.799.07"
$cell.ClearFormats()
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.506"
$cell.ClearFormats()
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "11.01"
$cell.ClearFormats()
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "2.315"
$cell.ClearFormats()
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "2.096.96"
$cell.ClearFormats()
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "159.56"
$cell.ClearFormats()
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "19.87"
$cell.ClearFormats()
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.168"
$cell.ClearFormats()
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "5.491"
$cell.ClearFormats()
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "122.24"
$cell.ClearFormats()
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "0.9895"
$cell.ClearFormats()
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.09540"
$cell.ClearFormats()
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "3.638"
$cell.ClearFormats()
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "5.355"
$cell.ClearFormats()
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "1.354"
$cell.ClearFormats()
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.06086"
$cell.ClearFormats()
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.02245"
$cell.ClearFormats()
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "8.333"
$cell.ClearFormats()
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.6010"
$cell.ClearFormats()
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.1896"
$cell.ClearFormats()
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "10.46"
$cell.ClearFormats()
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "1.248"
$cell.ClearFormats()
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.5686"
$cell.ClearFormats()
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "12.18"
$cell.ClearFormats()
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "1.944"
$cell.ClearFormats()
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "0.06783"
$cell.ClearFormats()
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "112.65"
$cell.ClearFormats()
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "3.039"
$cell.ClearFormats()

# Column E (Volume 1h) updates for all affected rows
$ws.Range("E2").Value = "  +0.16%  "
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -0.42%  "
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  -4.18%  "
$ws.Range("E10").Value = "  -0.52%  "
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("E14").Value = "  +0.84%  "
$ws.Range("E15").Value = "  +2.50%  "
$ws.Range("E16").Value = "  +0.01%  "
$ws.Range("E17").Value = "  +2.15%  "
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("E19").Value = "  +0.11%  "
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("E23").Value = "  -0.21%  "
$ws.Range("E24").Value = "  +0.09%  "
$ws.Range("E25").Value = "  -0.71%  "
$ws.Range("E26").Value = "  -0.54%  "
$ws.Range("E27").Value = "  +0.39%  "
$ws.Range("E28").Value = "  -1.47%  "
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("E30").Value = "  -1.32%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("E32").Value = "  +0.89%  "
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("E36").Value = "  -6.60%  "
$ws.Range("E37").Value = "  -1.10%  "
$ws.Range("E38").Value = "  -0.77%  "
$ws.Range("E39").Value = "  -1.34%  "
$ws.Range("E40").Value = "  +3.25%  "
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +0.16%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("E44").Value = "  +1.34%  "
$ws.Range("E45").Value = "  -0.91%  "
$ws.Range("E46").Value = "  -0.37%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("E49").Value = "  -1.85%  "
$ws.Range("E50").Value = "  -1.28%  "
$ws.Range("E51").Value = "  -10.60%  "
